# Adjust AffectorValue multipliers on the "AffectorValueLevelTable" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AffectorValueLevelTable")

$ws.Range("E3").Value = 0.5625
$ws.Range("E4").Value = 0.33333333300000001
$ws.Range("E5").Value = 1
